$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '274.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.67%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.69'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.72%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.865'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.71%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06321'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.89%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.890'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.46%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.280'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '36.06%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8700'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.99%'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1460'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.59%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05040'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.05%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07391'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.42%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02920'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-7.67%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09043'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.03%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001572'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.31%'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006304'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.64%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005986'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.40%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.448'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.09%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.323'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.87%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.296'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.45%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.10%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.37%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04368'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.28%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001175'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.04%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004261'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.30%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.04%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001693'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.20%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04043'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.31%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006657'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.96%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.10%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.79%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01218'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-12.92%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005307'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.51%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.488'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-37.09%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-33.07%'
